# Reporte_Regional.xlsx - actualizado 27 de abril de 2020
# Update the COVID-19 regional data table (B,C,D columns = nuevos/totales
# casos/fallecidos, H column = fecha) for every region row, then drop the
# stray formatted-but-empty rows (26, 27, 31, 36, 37) below the table so the
# sheet's used range shrinks back down to A1:H17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => Casos_nuevos(B), Casos_totales(C), Fallecidos(D), fecha(H)
$updates = @{
    2  = @(13,  265,  3, 43948)
    3  = @(8,   164,  1, 43948)
    4  = @(14,  457,  4, 43948)
    5  = @(6,   35,   0, 43948)
    6  = @(1,   74,   0, 43948)
    7  = @(25,  485,  9, 43948)
    8  = @(362, 7858, 95, 43948)
    9  = @(2,   94,   1, 43948)
    10 = @(9,   363,  12, 43948)
    11 = @(10,  741,  14, 43948)
    12 = @(3,   706,  6, 43948)
    13 = @(20,  1236, 32, 43948)
    14 = @(2,   180,  3, 43948)
    15 = @(4,   477,  8, 43948)
    16 = @(0,   7,    0, 43948)
    17 = @(3,   671,  10, 43948)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 8).Value = $vals[3]
}

# Remove the leftover stray cells (rows 26,27,31,36,37) that only carried
# formatting past the bottom of the real table -- drop the whole block so
# the sheet dimension goes back to A1:H17.
$ws.Range("A18:H37").EntireRow.Delete()
